$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'70.020.04"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.54%  "

$ws.Range("D3").Value = "'3.538.08"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.45%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").Value = "'604.18"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.40%  "

$ws.Range("D6").Value = "'196.72"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.42%  "

$ws.Range("D7").Value = "'0.627"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.88%  "

$ws.Range("E8").Value = "  -0.05%  "

$ws.Range("D9").Value = "'0.208"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.37%  "

$ws.Range("D10").Value = "'0.653"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.62%  "

$ws.Range("D11").Value = "'54.07"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.24%  "

$ws.Range("D12").Value = "'0.0000303"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.28%  "

$ws.Range("D13").Value = "'9.54"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.29%  "

$ws.Range("D14").Value = "'4.099.49"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.43%  "

$ws.Range("D15").Value = "'605.53"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.05%  "

$ws.Range("D16").Value = "'19.24"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.44%  "

$ws.Range("D17").Value = "'12.85"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.14%  "

$ws.Range("D18").Value = "'70.126.88"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.39%  "

$ws.Range("D19").Value = "'3.537.71"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.97%  "

$ws.Range("E20").Value = "  +0.31%  "

$ws.Range("D21").Value = "'0.997"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.63%  "

$ws.Range("D22").Value = "'18.07"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.87%  "

$ws.Range("E23").Value = "  +3.98%  "

$ws.Range("D24").Value = "'102.60"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.97%  "

$ws.Range("E25").Value = "  -2.44%  "

$ws.Range("D26").Value = "'3.12"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.41%  "

$ws.Range("D27").Value = "'10.95"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.65%  "

$ws.Range("D28").Value = "'9.63"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.96%  "

$ws.Range("D29").Value = "'33.64"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.30%  "

$ws.Range("D30").Value = "'7.15"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.46%  "

$ws.Range("D31").Value = "'4.34"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +15.57%  "

$ws.Range("D32").Value = "'12.68"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.74%  "

$ws.Range("E33").Value = "  -1.67%  "

$ws.Range("D34").Value = "'63.31"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.63%  "

$ws.Range("D35").Value = "'0.0₃0852"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +8.41%  "

$ws.Range("D36").Value = "'3.772.17"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +6.26%  "

$ws.Range("B37").Value = "Fetch.AI"
$ws.Range("C37").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D37").Value = "'3.08"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.75%  "

$ws.Range("B38").Value = "Dai"
$ws.Range("C38").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D38").Value = "'1.00"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.01%  "

$ws.Range("D39").Value = "'3.66"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.32%  "

$ws.Range("D40").Value = "'0.395"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.94%  "

$ws.Range("D41").Value = "'36.72"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.85%  "

$ws.Range("D42").Value = "'490.49"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -8.73%  "

$ws.Range("D43").Value = "'0.134"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.49%  "

$ws.Range("E44").Value = "  -2.89%  "

$ws.Range("E45").Value = "  -4.10%  "

$ws.Range("E46").Value = "  -3.11%  "

$ws.Range("D47").Value = "'3.30"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.22%  "

$ws.Range("E48").Value = "  +0.30%  "

$ws.Range("D49").Value = "'8.65"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.74%  "

$ws.Range("D50").Value = "'0.000250"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.36%  "

$ws.Range("D51").Value = "'130.36"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.85%  "
